$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: swap E3 and F3
$ws.Range("E3").Value = "Circuitos Elétricos"
$ws.Range("F3").Value = "-"

# Row 4: clear Circuitos Elétricos from E4
$ws.Range("E4").Value = "-"

# Row 6: move Desenho Técnico from C6 to D6, and set Circuitos Elétricos in E6
$ws.Range("C6").Value = "-"
$ws.Range("D6").Value = "Desenho Técnico"
$ws.Range("E6").Value = "Circuitos Elétricos"
